$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 482 }

for ($r = 2; $r -le $lastRow; $r++) {
    $pos = $ws.Cells.Item($r, 2).Value2
    if ($pos -eq "RB") {
        $cell = $ws.Cells.Item($r, 4)
        $cell.Value2 = $cell.Value2 + 3
    }
}
